# Amend the 4-phases pivot sheet so every currency-pair block also reports the
# 2015 interval (previously the earliest year shown was 2016/2017) and refreshes
# the count/percentage figures now that the extra year's occurrences are folded
# into the totals.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B to make room for the 2015 year column
$ws.Columns("B").Insert()
$ws.Columns("B").ColumnWidth = 5.2
# Columns H and I are now also year columns (shifted from the old total/percentage slots); narrow them to match the other year columns
$ws.Columns("H:I").ColumnWidth = 5.2

# Remove formatting artifacts left behind by the column insert on the currency-pair label rows
$ws.Range("B1").Clear()
$ws.Range("B8").Clear()
$ws.Range("B15").Clear()
$ws.Range("B22").Clear()

# ---- EURUSD ----
$ws.Range("B2").Value = 2015
$ws.Range("C2").Value = 2016
$ws.Range("D2").Value = 2017
$ws.Range("E2").Value = 2018
$ws.Range("F2").Value = 2019
$ws.Range("G2").Value = 2020
$ws.Range("H2").Value = 2021
$ws.Range("I2").Value = 2022
$ws.Range("J2").Value = "total_count_of_occurrences"
$ws.Range("K2").Value = "percentage_of_occurrences"
$ws.Range("B4").Value = 94
$ws.Range("C4").Value = 66
$ws.Range("D4").Value = 82
$ws.Range("E4").Value = 105
$ws.Range("F4").Value = 62
$ws.Range("G4").Value = 114
$ws.Range("H4").Value = 54
$ws.Range("I4").Value = 120
$ws.Range("J4").Value = 697
$ws.Range("K4").Value = 33.57418111753372
$ws.Range("B5").Value = 82
$ws.Range("C5").Value = 29
$ws.Range("D5").Value = 31
$ws.Range("E5").Value = 37
$ws.Range("F5").Value = 20
$ws.Range("G5").Value = 57
$ws.Range("H5").Value = 35
$ws.Range("I5").Value = 84
$ws.Range("J5").Value = 375
$ws.Range("K5").Value = 18.0635838150289
$ws.Range("B6").Value = 52
$ws.Range("C6").Value = 131
$ws.Range("D6").Value = 118
$ws.Range("E6").Value = 97
$ws.Range("F6").Value = 152
$ws.Range("G6").Value = 71
$ws.Range("H6").Value = 102
$ws.Range("I6").Value = 37
$ws.Range("J6").Value = 760
$ws.Range("K6").Value = 36.60886319845857
$ws.Range("B7").Value = 31
$ws.Range("C7").Value = 34
$ws.Range("D7").Value = 28
$ws.Range("E7").Value = 20
$ws.Range("F7").Value = 25
$ws.Range("G7").Value = 18
$ws.Range("H7").Value = 69
$ws.Range("I7").Value = 19
$ws.Range("J7").Value = 244
$ws.Range("K7").Value = 11.7533718689788

# ---- GBPUSD ----
$ws.Range("B9").Value = 2015
$ws.Range("C9").Value = 2016
$ws.Range("D9").Value = 2017
$ws.Range("E9").Value = 2018
$ws.Range("F9").Value = 2019
$ws.Range("G9").Value = 2020
$ws.Range("H9").Value = 2021
$ws.Range("I9").Value = 2022
$ws.Range("J9").Value = "total_count_of_occurrences"
$ws.Range("K9").Value = "percentage_of_occurrences"
$ws.Range("B11").Value = 108
$ws.Range("C11").Value = 125
$ws.Range("D11").Value = 77
$ws.Range("E11").Value = 83
$ws.Range("F11").Value = 52
$ws.Range("G11").Value = 109
$ws.Range("H11").Value = 52
$ws.Range("I11").Value = 136
$ws.Range("J11").Value = 742
$ws.Range("K11").Value = 35.74181117533718
$ws.Range("B12").Value = 50
$ws.Range("C12").Value = 36
$ws.Range("D12").Value = 21
$ws.Range("E12").Value = 57
$ws.Range("F12").Value = 61
$ws.Range("G12").Value = 56
$ws.Range("H12").Value = 21
$ws.Range("I12").Value = 64
$ws.Range("J12").Value = 366
$ws.Range("K12").Value = 17.63005780346821
$ws.Range("B13").Value = 86
$ws.Range("C13").Value = 87
$ws.Range("D13").Value = 133
$ws.Range("E13").Value = 86
$ws.Range("F13").Value = 95
$ws.Range("G13").Value = 75
$ws.Range("H13").Value = 160
$ws.Range("I13").Value = 50
$ws.Range("J13").Value = 772
$ws.Range("K13").Value = 37.1868978805395
$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 12
$ws.Range("D14").Value = 28
$ws.Range("E14").Value = 33
$ws.Range("F14").Value = 51
$ws.Range("G14").Value = 20
$ws.Range("H14").Value = 27
$ws.Range("I14").Value = 10
$ws.Range("J14").Value = 196
$ws.Range("K14").Value = 9.441233140655106

# ---- USDJPY ----
$ws.Range("B16").Value = 2015
$ws.Range("C16").Value = 2016
$ws.Range("D16").Value = 2017
$ws.Range("E16").Value = 2018
$ws.Range("F16").Value = 2019
$ws.Range("G16").Value = 2020
$ws.Range("H16").Value = 2021
$ws.Range("I16").Value = 2022
$ws.Range("J16").Value = "total_count_of_occurrences"
$ws.Range("K16").Value = "percentage_of_occurrences"
$ws.Range("B18").Value = 104
$ws.Range("C18").Value = 83
$ws.Range("D18").Value = 53
$ws.Range("E18").Value = 67
$ws.Range("F18").Value = 75
$ws.Range("G18").Value = 82
$ws.Range("H18").Value = 97
$ws.Range("I18").Value = 133
$ws.Range("J18").Value = 694
$ws.Range("K18").Value = 33.42967244701349
$ws.Range("B19").Value = 35
$ws.Range("C19").Value = 85
$ws.Range("D19").Value = 28
$ws.Range("E19").Value = 33
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 58
$ws.Range("H19").Value = 37
$ws.Range("I19").Value = 69
$ws.Range("J19").Value = 379
$ws.Range("K19").Value = 18.25626204238921
$ws.Range("B20").Value = 104
$ws.Range("C20").Value = 67
$ws.Range("D20").Value = 127
$ws.Range("E20").Value = 114
$ws.Range("F20").Value = 126
$ws.Range("G20").Value = 95
$ws.Range("H20").Value = 96
$ws.Range("I20").Value = 41
$ws.Range("J20").Value = 770
$ws.Range("K20").Value = 37.09055876685935
$ws.Range("B21").Value = 16
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 51
$ws.Range("E21").Value = 45
$ws.Range("F21").Value = 24
$ws.Range("G21").Value = 25
$ws.Range("H21").Value = 30
$ws.Range("I21").Value = 17
$ws.Range("J21").Value = 233
$ws.Range("K21").Value = 11.22350674373796

# ---- XAUUSD ----
$ws.Range("B23").Value = 2015
$ws.Range("C23").Value = 2016
$ws.Range("D23").Value = 2017
$ws.Range("E23").Value = 2018
$ws.Range("F23").Value = 2019
$ws.Range("G23").Value = 2020
$ws.Range("H23").Value = 2021
$ws.Range("I23").Value = 2022
$ws.Range("J23").Value = "total_count_of_occurrences"
$ws.Range("K23").Value = "percentage_of_occurrences"
$ws.Range("B25").Value = 109
$ws.Range("C25").Value = 85
$ws.Range("D25").Value = 74
$ws.Range("E25").Value = 102
$ws.Range("F25").Value = 94
$ws.Range("G25").Value = 122
$ws.Range("H25").Value = 67
$ws.Range("I25").Value = 93
$ws.Range("J25").Value = 746
$ws.Range("K25").Value = 36.14341085271317
$ws.Range("B26").Value = 26
$ws.Range("C26").Value = 41
$ws.Range("D26").Value = 21
$ws.Range("E26").Value = 19
$ws.Range("F26").Value = 50
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = 39
$ws.Range("I26").Value = 55
$ws.Range("J26").Value = 289
$ws.Range("K26").Value = 14.00193798449612
$ws.Range("B27").Value = 105
$ws.Range("C27").Value = 94
$ws.Range("D27").Value = 113
$ws.Range("E27").Value = 115
$ws.Range("F27").Value = 80
$ws.Range("G27").Value = 82
$ws.Range("H27").Value = 103
$ws.Range("I27").Value = 70
$ws.Range("J27").Value = 762
$ws.Range("K27").Value = 36.91860465116279
$ws.Range("B28").Value = 18
$ws.Range("C28").Value = 38
$ws.Range("D28").Value = 49
$ws.Range("E28").Value = 22
$ws.Range("F28").Value = 34
$ws.Range("G28").Value = 17
$ws.Range("H28").Value = 49
$ws.Range("I28").Value = 40
$ws.Range("J28").Value = 267
$ws.Range("K28").Value = 12.93604651162791
